$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: Status columns E & F show "Ready for handoff" -> new status
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.9777050018311
$overview.Columns.Item(6).ColumnWidth = 29.9777050018311

# zh-cn sheet
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("L2").Value = "2017-02-21 03:01:57"
$zhcn.Range("L3").Value = "2017-02-21 03:01:57"
$zhcn.Range("M2").Value = ""
$zhcn.Range("R2").Value = ""
$zhcn.Range("R3").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.9777050018311
$zhcn.Columns.Item(13).ColumnWidth = 23.8743762969971
$zhcn.Columns.Item(18).ColumnWidth = 13.7470531463623

# de-de sheet
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("L2").Value = "2017-02-21 03:02:19"
$dede.Range("L3").Value = "2017-02-21 03:02:19"
$dede.Range("M2").Value = ""
$dede.Range("R2").Value = ""
$dede.Range("R3").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.9777050018311
$dede.Columns.Item(13).ColumnWidth = 23.8743762969971
$dede.Columns.Item(18).ColumnWidth = 13.7470531463623
